$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing column B values (rows 2, 3, 6 change text; rows 4 & 5 stay the same)
$ws.Range("B2").Value = "長出"
$ws.Range("B3").Value = "190 gold"
$ws.Range("B6").Value = "battle_text longest_name"

# Add the new column C: header + per-row tag names
$ws.Range("C1").Value = "TextID"
$ws.Range("C2").Value = "年会"
$ws.Range("C3").Value = "SwordOfTruth_Price"
$ws.Range("C4").Value = "LongestName_TagTest"
$ws.Range("C5").Value = "BattleText_TagTest"
$ws.Range("C6").Value = "Multiple_TagTest"

# Update column widths so C matches the authored width, consistent with A/B bestFit cols
$ws.Columns.Item(3).ColumnWidth = 18.5546875

# Match the post-edit selection / active cell noted in the workbook view
$ws.Range("B12").Select() | Out-Null
